$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.135.42"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").Value = "1.572.52"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'207.24"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "'0.491"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'22.25"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").Value = "'0.0868"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "1.796.47"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "1.562.33"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "'3.77"
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").Value = "'0.517"
$ws.Range("E15").Value = "  -1.50%  "
$ws.Range("D16").Value = "27.148.86"
$ws.Range("E16").Value = "  -1.48%  "
$ws.Range("D17").Value = "'62.18"
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").Value = "'213.97"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "0.0₃0683"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").Value = "'9.42"
$ws.Range("E23").Value = "  -3.64%  "
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").Value = "'152.51"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("E26").Value = "  -3.34%  "
$ws.Range("D27").Value = "'14.92"
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("E30").Value = "  -3.23%  "
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("D32").Value = "'3.17"
$ws.Range("E32").Value = "  -1.60%  "
$ws.Range("D33").Value = "1.394.63"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("E37").Value = "  -2.92%  "
$ws.Range("E38").Value = "  -2.35%  "
$ws.Range("D39").Value = "'0.813"
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("E40").Value = "  -3.41%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("E42").Value = "  +3.82%  "
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("D44").Value = "'5.42"
$ws.Range("E44").Value = "  +2.22%  "
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").Value = "'63.68"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("D47").Value = "1.708.62"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").Value = "'85.49"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").Value = "0.0₇0990"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0499"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0951"
$ws.Range("E51").Value = "  -0.89%  "
